$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1850220264317181
$ws.Range("C2").Value = 0.6013215859030837
$ws.Range("J2").Value = 0.01541850220264317
$ws.Range("P2").Value = 0.1299559471365639
$ws.Range("S2").Value = 0.06828193832599119
$ws.Range("B3").Value = 0.003546099290780142
$ws.Range("C3").Value = 0.02127659574468085
$ws.Range("J3").Value = 0.02482269503546099
$ws.Range("P3").Value = 0.7624113475177305
$ws.Range("S3").Value = 0.1879432624113475
$ws.Range("J4").Value = 0.09615384615384616
$ws.Range("P4").Value = 0.6923076923076923
$ws.Range("S4").Value = 0.2115384615384615
$ws.Range("O5").Value = 0.5
$ws.Range("P5").Value = 0.5
$ws.Range("B6").Value = 0.06639004149377593
$ws.Range("D6").Value = 0.02074688796680498
$ws.Range("E6").Value = 0.008298755186721992
$ws.Range("F6").Value = 0.05809128630705394
$ws.Range("J6").Value = 0.2489626556016598
$ws.Range("O6").Value = 0.01244813278008299
$ws.Range("Q6").Value = 0.1493775933609958
$ws.Range("R6").Value = 0.08713692946058091
$ws.Range("S6").Value = 0.3485477178423236
$ws.Range("B7").Value = 0.1448763250883392
$ws.Range("D7").Value = 0.03180212014134275
$ws.Range("F7").Value = 0.08480565371024736
$ws.Range("J7").Value = 0.1060070671378092
$ws.Range("O7").Value = 0.01060070671378092
$ws.Range("Q7").Value = 0.1484098939929329
$ws.Range("R7").Value = 0.06360424028268551
$ws.Range("S7").Value = 0.4098939929328622
$ws.Range("B8").Value = 0.1255144032921811
$ws.Range("D8").Value = 0.01851851851851852
$ws.Range("F8").Value = 0.06172839506172839
$ws.Range("J8").Value = 0.1234567901234568
$ws.Range("O8").Value = 0.01851851851851852
$ws.Range("Q8").Value = 0.1460905349794239
$ws.Range("R8").Value = 0.08847736625514403
$ws.Range("S8").Value = 0.4176954732510288
$ws.Range("B9").Value = 0.1140350877192982
$ws.Range("D9").Value = 0.02631578947368421
$ws.Range("F9").Value = 0.07017543859649122
$ws.Range("J9").Value = 0.1754385964912281
$ws.Range("O9").Value = 0.05263157894736842
$ws.Range("Q9").Value = 0.1403508771929824
$ws.Range("R9").Value = 0.07017543859649122
$ws.Range("S9").Value = 0.3508771929824561
$ws.Range("B10").Value = 0.150032195750161
$ws.Range("D10").Value = 0.0180296200901481
$ws.Range("F10").Value = 0.05988409529942048
$ws.Range("J10").Value = 0.1313586606567933
$ws.Range("O10").Value = 0.0154539600772698
$ws.Range("Q10").Value = 0.1854475209272376
$ws.Range("R10").Value = 0.08564069542820348
$ws.Range("S10").Value = 0.3541532517707663
$ws.Range("G11").Value = 0.1352549889135255
$ws.Range("J11").Value = 0.08425720620842572
$ws.Range("K11").Value = 0.1951219512195122
$ws.Range("L11").Value = 0.5676274944567627
$ws.Range("S11").Value = 0.01773835920177384
$ws.Range("G12").Value = 0.7213740458015268
$ws.Range("J12").Value = 0.2022900763358779
$ws.Range("K12").Value = 0.01145038167938931
$ws.Range("L12").Value = 0.01526717557251908
$ws.Range("S12").Value = 0.04961832061068702
$ws.Range("G13").Value = 0.6461538461538462
$ws.Range("J13").Value = 0.2615384615384616
$ws.Range("S13").Value = 0.09230769230769231
$ws.Range("F15").Value = 0.004694835680751174
$ws.Range("H15").Value = 0.1408450704225352
$ws.Range("I15").Value = 0.04694835680751173
$ws.Range("J15").Value = 0.3708920187793427
$ws.Range("K15").Value = 0.03286384976525822
$ws.Range("M15").Value = 0.01408450704225352
$ws.Range("N15").Value = 0.004694835680751174
$ws.Range("O15").Value = 0.04694835680751173
$ws.Range("S15").Value = 0.3380281690140845
$ws.Range("F16").Value = 0.0198019801980198
$ws.Range("H16").Value = 0.2145214521452145
$ws.Range("I16").Value = 0.05280528052805281
$ws.Range("J16").Value = 0.4026402640264026
$ws.Range("K16").Value = 0.132013201320132
$ws.Range("M16").Value = 0.0132013201320132
$ws.Range("O16").Value = 0.03630363036303631
$ws.Range("S16").Value = 0.1287128712871287
$ws.Range("F17").Value = 0.01339285714285714
$ws.Range("H17").Value = 0.1674107142857143
$ws.Range("I17").Value = 0.04464285714285714
$ws.Range("J17").Value = 0.4107142857142857
$ws.Range("K17").Value = 0.1584821428571428
$ws.Range("M17").Value = 0.02232142857142857
$ws.Range("O17").Value = 0.06026785714285714
$ws.Range("S17").Value = 0.1227678571428571
$ws.Range("F18").Value = 0.01801801801801802
$ws.Range("H18").Value = 0.1486486486486487
$ws.Range("I18").Value = 0.05855855855855856
$ws.Range("J18").Value = 0.4189189189189189
$ws.Range("K18").Value = 0.1486486486486487
$ws.Range("M18").Value = 0.03153153153153153
$ws.Range("O18").Value = 0.07207207207207207
$ws.Range("S18").Value = 0.1036036036036036
$ws.Range("F19").Value = 0.01512027491408935
$ws.Range("H19").Value = 0.1965635738831615
$ws.Range("I19").Value = 0.03711340206185567
$ws.Range("J19").Value = 0.4075601374570447
$ws.Range("K19").Value = 0.1415807560137457
$ws.Range("M19").Value = 0.0288659793814433
$ws.Range("N19").Value = 0.001374570446735395
$ws.Range("O19").Value = 0.05429553264604811
$ws.Range("S19").Value = 0.1175257731958763
